$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") for rows 2 through 27 all currently hold the
# serial date value 45233 (2023-11-03) and need to be updated to 45243
# (2023-11-13).
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
